# Applies a scheduled market-data refresh to the cached price/profit
# columns (H:N) on several sheets. Values come from a fresh Universalis
# pull, so most cells are simple overwrites; ALC!N101 is no longer
# populated by the refreshed data (cell cleared), and LTW!M136, WVR!M15
# and WVR!M132 newly gain a cached profit figure.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 9773.0625
$ws.Range("I80").Value = 346.77777
$ws.Range("J80").Value = 21892.572
$ws.Range("K80").Value = 1040.33331
$ws.Range("L80").Value = 65677.716
$ws.Range("M80").Value = -42.33330999999998
$ws.Range("N80").Value = -67673.716
$ws.Range("H83").Value = 9773.0625
$ws.Range("I83").Value = 346.77777
$ws.Range("J83").Value = 21892.572
$ws.Range("K83").Value = 3120.99993
$ws.Range("L83").Value = 197033.148
$ws.Range("M83").Value = 1871.00007
$ws.Range("N83").Value = -207017.148
$ws.Range("H101").Value = 408
$ws.Range("I101").Value = 408
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 1224
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 398
$ws.Range("N101").ClearContents()  # removed (was -5642.5)
$ws.Range("H112").Value = 2381
$ws.Range("J112").Value = 2727.4707
$ws.Range("L112").Value = 8182.4121
$ws.Range("N112").Value = -10398.4121
$ws.Range("H116").Value = 8913.25
$ws.Range("I116").Value = 8884.666999999999
$ws.Range("K116").Value = 8884.666999999999
$ws.Range("M116").Value = -5442.666999999999
$ws.Range("H132").Value = 1415.3429
$ws.Range("I132").Value = 1368.2333
$ws.Range("K132").Value = 4104.699900000001
$ws.Range("M132").Value = -1574.699900000001
$ws.Range("H137").Value = 2645.3333
$ws.Range("I137").Value = 1830.0588
$ws.Range("K137").Value = 5490.1764
$ws.Range("M137").Value = -2940.1764
$ws.Range("H138").Value = 6952.12
$ws.Range("I138").Value = 1545.6666
$ws.Range("J138").Value = 9993.25
$ws.Range("K138").Value = 4636.9998
$ws.Range("L138").Value = 29979.75
$ws.Range("M138").Value = 503.0002000000004
$ws.Range("N138").Value = -40259.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1397.2222
$ws.Range("I61").Value = 1397.2222
$ws.Range("K61").Value = 1397.2222
$ws.Range("M61").Value = -1185.2222
$ws.Range("H132").Value = 1607.1818
$ws.Range("I132").Value = 1534.9672
$ws.Range("K132").Value = 4604.9016
$ws.Range("M132").Value = -2074.9016
$ws.Range("H136").Value = 1397.2222
$ws.Range("I136").Value = 1397.2222
$ws.Range("K136").Value = 4191.6666
$ws.Range("M136").Value = -1641.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2999
$ws.Range("I86").Value = 2999
$ws.Range("K86").Value = 2999
$ws.Range("M86").Value = -1876
$ws.Range("H89").Value = 2999
$ws.Range("I89").Value = 2999
$ws.Range("K89").Value = 14995
$ws.Range("M89").Value = -9379
$ws.Range("H105").Value = 3635.658
$ws.Range("I105").Value = 3006.682
$ws.Range("J105").Value = 4500.5
$ws.Range("K105").Value = 3006.682
$ws.Range("L105").Value = 4500.5
$ws.Range("M105").Value = -1259.682
$ws.Range("N105").Value = -7994.5
$ws.Range("H107").Value = 1388.375
$ws.Range("I107").Value = 962.6923
$ws.Range("J107").Value = 3233
$ws.Range("K107").Value = 962.6923
$ws.Range("L107").Value = 3233
$ws.Range("M107").Value = 957.3077
$ws.Range("N107").Value = -7073
$ws.Range("H134").Value = 597.6
$ws.Range("I134").Value = 586.0417
$ws.Range("J134").Value = 875
$ws.Range("K134").Value = 1758.1251
$ws.Range("L134").Value = 2625
$ws.Range("M134").Value = 776.8749
$ws.Range("N134").Value = -7695

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 3335236.8
$ws.Range("I6").Value = 5004183
$ws.Range("K6").Value = 5004183
$ws.Range("M6").Value = -5004070
$ws.Range("H31").Value = 4877.4443
$ws.Range("I31").Value = 2989.3333
$ws.Range("K31").Value = 2989.3333
$ws.Range("M31").Value = -2694.3333
$ws.Range("H34").Value = 4877.4443
$ws.Range("I34").Value = 2989.3333
$ws.Range("K34").Value = 2989.3333
$ws.Range("M34").Value = -2787.3333
$ws.Range("H64").Value = 24000
$ws.Range("J64").Value = 24000
$ws.Range("L64").Value = 24000
$ws.Range("N64").Value = -24496
$ws.Range("H67").Value = 24000
$ws.Range("J67").Value = 24000
$ws.Range("L67").Value = 24000
$ws.Range("N67").Value = -25716
$ws.Range("H132").Value = 2467.925
$ws.Range("I132").Value = 2115.3
$ws.Range("K132").Value = 6345.900000000001
$ws.Range("M132").Value = -3815.900000000001
$ws.Range("H134").Value = 2654.3333
$ws.Range("I134").Value = 1701.7916
$ws.Range("K134").Value = 5105.3748
$ws.Range("M134").Value = -2570.3748

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2028.8889
$ws.Range("J131").Value = 2028.8889
$ws.Range("L131").Value = 6086.6667
$ws.Range("N131").Value = -16166.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 276.125
$ws.Range("J2").Value = 408.5
$ws.Range("L2").Value = 408.5
$ws.Range("N2").Value = -634.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5287
$ws.Range("I16").Value = 5287
$ws.Range("K16").Value = 5287
$ws.Range("M16").Value = -5117
$ws.Range("H40").Value = 3832.3333
$ws.Range("I40").Value = 3832.3333
$ws.Range("K40").Value = 3832.3333
$ws.Range("M40").Value = -3696.3333
$ws.Range("H68").Value = 3387.25
$ws.Range("I68").Value = 3033
$ws.Range("K68").Value = 3033
$ws.Range("M68").Value = -2284
$ws.Range("H71").Value = 3387.25
$ws.Range("I71").Value = 3033
$ws.Range("K71").Value = 15165
$ws.Range("M71").Value = -11421
$ws.Range("H132").Value = 4160.4165
$ws.Range("I132").Value = 2703.5715
$ws.Range("K132").Value = 8110.7145
$ws.Range("M132").Value = -5580.7145
$ws.Range("H136").Value = 3027
$ws.Range("I136").Value = 3027
$ws.Range("K136").Value = 9081
$ws.Range("M136").Value = -6531

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 50000
$ws.Range("I15").Value = 50000
$ws.Range("K15").Value = 50000
$ws.Range("M15").Value = -49712
$ws.Range("H132").Value = 3250
$ws.Range("I132").Value = 2500
$ws.Range("K132").Value = 7500
$ws.Range("M132").Value = -4970
$ws.Range("H136").Value = 38231.93
$ws.Range("I136").Value = 2205.05
$ws.Range("J136").Value = 128299.125
$ws.Range("K136").Value = 6615.150000000001
$ws.Range("L136").Value = 384897.375
$ws.Range("M136").Value = -4065.150000000001
$ws.Range("N136").Value = -389997.375
